$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number / week-covering dates) ---
$ws.Range("A8").Value = "Volume 30   Number  40"
$ws.Range("C9").Value = "Report Covering the Week  10/2/2023  Through  10/8/2023"

# --- Weekly crime-stat numeric updates (rows 14-30) ---
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = 0
$ws.Range("G14").Value = 11
$ws.Range("H14").Value = -45.454545454545
$ws.Range("I14").Value = 56
$ws.Range("J14").Value = 61
$ws.Range("K14").Value = -8.196721311475
$ws.Range("L14").Value = -28.205128205128
$ws.Range("M14").Value = -50.442477876106
$ws.Range("N14").Value = -85.340314136125
$ws.Range("D15").Value = 7
$ws.Range("E15").Value = -42.857142857142
$ws.Range("F15").Value = 17
$ws.Range("G15").Value = 29
$ws.Range("H15").Value = -41.379310344827
$ws.Range("I15").Value = 173
$ws.Range("J15").Value = 200
$ws.Range("K15").Value = -13.5
$ws.Range("L15").Value = 2.366863905325
$ws.Range("M15").Value = 6.134969325153
$ws.Range("N15").Value = -63.269639065817
$ws.Range("C16").Value = 51
$ws.Range("D16").Value = 53
$ws.Range("E16").Value = -3.77358490566
$ws.Range("F16").Value = 189
$ws.Range("G16").Value = 199
$ws.Range("H16").Value = -5.02512562814
$ws.Range("I16").Value = 1883
$ws.Range("J16").Value = 2016
$ws.Range("K16").Value = -6.597222222222
$ws.Range("L16").Value = 21.483870967741
$ws.Range("M16").Value = -30.797500918779
$ws.Range("N16").Value = -85.200031439126
$ws.Range("C17").Value = 93
$ws.Range("D17").Value = 68
$ws.Range("E17").Value = 36.764705882352
$ws.Range("F17").Value = 355
$ws.Range("G17").Value = 296
$ws.Range("H17").Value = 19.932432432432
$ws.Range("I17").Value = 3321
$ws.Range("J17").Value = 3234
$ws.Range("K17").Value = 2.690166975881
$ws.Range("L17").Value = 20.282506338283
$ws.Range("M17").Value = 26.707363601678
$ws.Range("N17").Value = -50.528824668553
$ws.Range("C18").Value = 44
$ws.Range("D18").Value = 49
$ws.Range("E18").Value = -10.204081632653
$ws.Range("F18").Value = 162
$ws.Range("G18").Value = 190
$ws.Range("H18").Value = -14.736842105263
$ws.Range("I18").Value = 1590
$ws.Range("J18").Value = 1849
$ws.Range("K18").Value = -14.007571660357
$ws.Range("L18").Value = 3.046014257939
$ws.Range("M18").Value = -35.286935286935
$ws.Range("N18").Value = -83.049040511727
$ws.Range("C19").Value = 121
$ws.Range("D19").Value = 119
$ws.Range("E19").Value = 1.680672268907
$ws.Range("F19").Value = 447
$ws.Range("G19").Value = 478
$ws.Range("H19").Value = -6.485355648535
$ws.Range("I19").Value = 4455
$ws.Range("J19").Value = 4583
$ws.Range("K19").Value = -2.792930394937
$ws.Range("L19").Value = 27.322092026293
$ws.Range("M19").Value = 36.908420405654
$ws.Range("N19").Value = -16.573033707865
$ws.Range("C20").Value = 33
$ws.Range("D20").Value = 32
$ws.Range("E20").Value = 3.125
$ws.Range("F20").Value = 165
$ws.Range("G20").Value = 166
$ws.Range("H20").Value = -0.602409638554
$ws.Range("I20").Value = 1416
$ws.Range("J20").Value = 1419
$ws.Range("K20").Value = -0.211416490486
$ws.Range("L20").Value = 20.716112531969
$ws.Range("M20").Value = 28.260869565217
$ws.Range("N20").Value = -80.455486542443
$ws.Range("C21").Value = 348
$ws.Range("D21").Value = 330
$ws.Range("E21").Value = 5.454545454545
$ws.Range("F21").Value = 1341
$ws.Range("G21").Value = 1369
$ws.Range("H21").Value = -2.045288531775
$ws.Range("I21").Value = 12894
$ws.Range("J21").Value = 13362
$ws.Range("K21").Value = -3.502469690166
$ws.Range("L21").Value = 19.688109161793
$ws.Range("M21").Value = 3.707874205742
$ws.Range("N21").Value = -69.484545841813
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 6
$ws.Range("E22").Value = -66.666666666666
$ws.Range("F22").Value = 20
$ws.Range("G22").Value = 23
$ws.Range("H22").Value = -13.043478260869
$ws.Range("I22").Value = 218
$ws.Range("J22").Value = 269
$ws.Range("K22").Value = -18.959107806691
$ws.Range("L22").Value = 11.794871794871
$ws.Range("M22").Value = -33.333333333333
$ws.Range("C23").Value = 30
$ws.Range("D23").Value = 27
$ws.Range("E23").Value = 11.111111111111
$ws.Range("F23").Value = 114
$ws.Range("G23").Value = 105
$ws.Range("H23").Value = 8.571428571428
$ws.Range("I23").Value = 1215
$ws.Range("J23").Value = 1184
$ws.Range("K23").Value = 2.618243243243
$ws.Range("L23").Value = 8.482142857142
$ws.Range("M23").Value = 33.369923161361
$ws.Range("C24").Value = 205
$ws.Range("D24").Value = 281
$ws.Range("E24").Value = -27.046263345195
$ws.Range("F24").Value = 941
$ws.Range("G24").Value = 1178
$ws.Range("H24").Value = -20.118845500848
$ws.Range("I24").Value = 9621
$ws.Range("J24").Value = 10351
$ws.Range("K24").Value = -7.052458699642
$ws.Range("L24").Value = 23.346153846153
$ws.Range("M24").Value = 19.977553310886
$ws.Range("C25").Value = 118
$ws.Range("D25").Value = 74
$ws.Range("E25").Value = 59.459459459459
$ws.Range("F25").Value = 458
$ws.Range("G25").Value = 415
$ws.Range("H25").Value = 10.361445783132
$ws.Range("I25").Value = 4754
$ws.Range("J25").Value = 4578
$ws.Range("K25").Value = 3.844473569244
$ws.Range("L25").Value = 32.756213348226
$ws.Range("M25").Value = -23.223514211886
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -40
$ws.Range("F26").Value = 26
$ws.Range("G26").Value = 39
$ws.Range("H26").Value = -33.333333333333
$ws.Range("I26").Value = 266
$ws.Range("J26").Value = 298
$ws.Range("K26").Value = -10.738255033557
$ws.Range("L26").Value = -9.523809523809
$ws.Range("C27").Value = 16
$ws.Range("D27").Value = 12
$ws.Range("E27").Value = 33.333333333333
$ws.Range("F27").Value = 64
$ws.Range("G27").Value = 40
$ws.Range("H27").Value = 60
$ws.Range("I27").Value = 503
$ws.Range("J27").Value = 474
$ws.Range("K27").Value = 6.118143459915
$ws.Range("L27").Value = -6.33147113594
$ws.Range("C28").Value = 7
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 250
$ws.Range("F28").Value = 19
$ws.Range("G28").Value = 21
$ws.Range("H28").Value = -9.523809523809
$ws.Range("I28").Value = 188
$ws.Range("J28").Value = 276
$ws.Range("K28").Value = -31.884057971014
$ws.Range("L28").Value = -45.029239766081
$ws.Range("M28").Value = -55.971896955503
$ws.Range("N28").Value = -87.647831800262
$ws.Range("C29").Value = 5
$ws.Range("D29").Value = 2
$ws.Range("E29").Value = 150
$ws.Range("F29").Value = 14
$ws.Range("G29").Value = 16
$ws.Range("H29").Value = -12.5
$ws.Range("I29").Value = 158
$ws.Range("J29").Value = 229
$ws.Range("K29").Value = -31.004366812227
$ws.Range("L29").Value = -42.335766423357
$ws.Range("M29").Value = -54.335260115606
$ws.Range("N29").Value = -88.467153284671
$ws.Range("G30").Value = 12
$ws.Range("H30").Value = -66.666666666666
$ws.Range("I30").Value = 49
$ws.Range("K30").Value = -26.865671641791
$ws.Range("L30").Value = 0

# --- Row 30 (Other Violations) now reports suppressed week-to-date figures as "0" / "N/A" text,
#     matching the style used elsewhere in the sheet for text placeholders (e.g. M30/N30).
$ws.Range("C30").Value = "'0"
$ws.Range("D30").Value = "'0"
$ws.Range("E30").Value = "'***.*"
$ws.Range("M30").Copy()
$ws.Range("C30:E30").PasteSpecial(-4122)
